$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Asgardians Slot for Free - A
#    Review of the Game" right before the final "Prompt: ..." paragraph.
# ------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($lastIdx)
$promptPara.Range.InsertParagraphBefore()

$newIdx = $lastIdx
$newPara = $d.Paragraphs.Item($newIdx)
$newRange = $newPara.Range

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Asgardians Slot for Free - A Review of the Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newRange.InsertXML($xmlFrag)

# ------------------------------------------------------------------
# 3) Replace the old image-prompt text of the final paragraph with the
#    meta-description sentence (formatting/italics are preserved).
# ------------------------------------------------------------------
$oldPromptText = 'Prompt: Create a feature image for the online slot game "Asgardians" that features a happy Maya warrior with glasses in cartoon style. The image should showcase the adventurous and exciting nature of the game, as well as its Norse mythology theme. Please use bright colors and dynamic visual elements to catch viewers'' attention and encourage them to play the game. The image should also include the title "Asgardians" in bold font to clearly convey the game''s identity.'
$newPromptText = 'Read our review of the Asgardians slot game. Play for free and win up to 7,500x your bet in a single spin with this high-paying game inspired by Norse mythology.'

$null = $d.Content.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, $newPromptText, 2)
